# Apply the edits described by the commit:
#  - Rename sheet "Feuil1" -> "TIPE_Homme_Badminton_Simple"
#  - Update some tournament-stage labels in column A (rows 13-21) from
#    "Poules préliminaires" to "Finale" / "Demi finale" / "Quart de finale"
#  - Move the active cell selection to F35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "TIPE_Homme_Badminton_Simple"

# Update column A labels for rows 13-21
$ws.Range("A13").Value = "Finale"
$ws.Range("A14").Value = "Demi finale"
$ws.Range("A15").Value = "Quart de finale"
$ws.Range("A16").Value = "Demi finale"
$ws.Range("A17").Value = "Quart de finale"
$ws.Range("A18").Value = "Finale"
$ws.Range("A19").Value = "Demi finale"
$ws.Range("A20").Value = "Quart de finale"
$ws.Range("A21").Value = "Demi finale"

# Move the selection to F35 to match the saved view state
$ws.Range("F35").Select()
